# "letras pablo con etiqueta"
# Update the label/value column (AR, "valor") for every data row from 25 to 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 51; $r++) {
    $ws.Range("AR$r").Value = 15
}

# Reflect the view state captured in the saved file (scrolled/selected range).
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Application.ActiveWindow.ScrollColumn = 37
$ws.Range("AW35").Select()
